$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the repayment figures on row 3 (Original / Outstanding / Over Due)
$ws.Range("A3").Value = 256.26
$ws.Range("E3").Value = 166.98
$ws.Range("F3").Value = 166.98

# Move the active selection to D5, matching the saved view state
$ws.Range("D5").Select()
